$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new annotation columns to the header row (row 1)
$ws.Range("I1").Value = "annotations.de"
$ws.Range("J1").Value = "annotations.en"

# Move the selection to K1, matching the post-edit cursor position
$ws.Range("K1").Select()
